$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the Gantt-chart style row: a new schedule entry at row 11
$ws.Range("A11").Value = Get-Date -Year 2018 -Month 10 -Day 8 -Hour 0 -Minute 0 -Second 0
$ws.Range("B11").Value = "2:00PM"
$ws.Range("C11").Value = "Student Union"

# Header cell A5 gets centered alignment (new style)
$ws.Range("A5").HorizontalAlignment = -4108

# Update the view: scroll position, zoom level, and active selection
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 150
$ws.Range("C14").Select() | Out-Null
